# Reprocess metadata sheet with newly curated dimensions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-* type row. "situacion-preferente" and "sexo" become measures
# instead of dimensions; "aragon" becomes a sdmx-dimension:refArea just like
# comarca/provincia.
$ws.Range("B2").Value2 = "iaest-measure:situacion-preferente"
$ws.Range("E2").Value2 = "sdmx-dimension:refArea"
$ws.Range("F2").Value2 = "iaest-measure:sexo"

# Row 3: medida/dim classification row, follows the row 2 re-classification.
$ws.Range("A3").Value2 = "medida"
$ws.Range("B3").Value2 = "medida"
$ws.Range("D3").Value2 = "dim"
$ws.Range("E3").Value2 = "dim"
$ws.Range("F3").Value2 = "medida"
$ws.Range("H3").Value2 = "medida"
$ws.Range("I3").Value2 = "dim"

# Row 4: value-type / URI row, follows the row 2/3 re-classification.
$ws.Range("A4").Value2 = "xsd:int"
$ws.Range("B4").Value2 = "xsd:int"
$ws.Range("D4").Value2 = "URI-Provincia"
$ws.Range("E4").Value2 = "URI-Comunidad"
$ws.Range("F4").Value2 = "xsd:int"
$ws.Range("H4").Value2 = "xsd:int"
$ws.Range("I4").Value2 = "URI-comarca"

# Row 5 (mapping-*.xlsx references) is no longer needed; remove it entirely.
$ws.Rows.Item(5).Delete()
